$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-06-17 Tuesday" "2025-06-18 Wednesday"

Replace-Text "43÷2=21, 1" "62÷9=6, 8"
Replace-Text "89÷6=14, 5" "25÷9=2, 7"
Replace-Text "12÷6=2, 0" "60÷6=10, 0"
Replace-Text "41÷2=20, 1" "79÷8=9, 7"
Replace-Text "57÷7=8, 1" "35÷7=5, 0"

Replace-Text "99÷8=12, 3" "67÷3=22, 1"
Replace-Text "36÷2=18, 0" "88÷5=17, 3"
Replace-Text "75÷2=37, 1" "18÷2=9, 0"
Replace-Text "97÷6=16, 1" "41÷9=4, 5"
Replace-Text "56÷4=14, 0" "23÷9=2, 5"

Replace-Text "77÷9=8, 5" "54÷8=6, 6"
Replace-Text "63÷6=10, 3" "71÷5=14, 1"
Replace-Text "34÷9=3, 7" "11÷5=2, 1"
Replace-Text "99÷9=11, 0" "96÷6=16, 0"
Replace-Text "19÷9=2, 1" "74÷8=9, 2"

Replace-Text "83÷9=9, 2" "66÷6=11, 0"
Replace-Text "15÷9=1, 6" "36÷7=5, 1"
Replace-Text "70÷9=7, 7" "90÷6=15, 0"
Replace-Text "68÷9=7, 5" "40÷2=20, 0"
Replace-Text "10÷7=1, 3" "12÷2=6, 0"

Replace-Text "29÷7=4, 1" "83÷8=10, 3"
Replace-Text "10÷2=5, 0" "76÷9=8, 4"
Replace-Text "28÷2=14, 0" "29÷5=5, 4"
Replace-Text "37÷6=6, 1" "87÷8=10, 7"
Replace-Text "24÷8=3, 0" "30÷6=5, 0"

Write-Output "Done"
